# Duplicate "leafId" column was found in the data - remove it.
# Column A ("leafId") duplicated the info already represented by the
# "vendorLeaf" column, so delete the entire column A, shifting the
# remaining columns (vendorLeaf, vendor, purchasedDate) left.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the whole column A (leafId) - this shifts B->A, C->B, D->C
$ws.Range("A:A").EntireColumn.Delete()

# Update the selected/active cell as recorded after the edit
$ws.Range("D10").Select()
